$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.364.24"
$ws.Range("E2").Value = "  -1.09%  "

$ws.Range("D3").Value = "2.514.66"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.51"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.02"
$ws.Range("E6").Value = "  -1.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.21"
$ws.Range("E10").Value = "  -1.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.64"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("E13").Value = "  -3.25%  "

$ws.Range("D14").Value = "2.898.75"
$ws.Range("E14").Value = "  -2.31%  "

$ws.Range("D15").Value = "2.547.60"
$ws.Range("E15").Value = "  -3.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.21"
$ws.Range("E16").Value = "  +5.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("E17").Value = "  -2.57%  "

$ws.Range("D18").Value = "42.314.88"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.92"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").Value = "0.0₃0972"
$ws.Range("E20").Value = "  -1.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.42"
$ws.Range("E21").Value = "  -2.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.27"
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.89"
$ws.Range("E23").Value = "  -2.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.90"
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.02"
$ws.Range("E25").Value = "  -4.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.16"
$ws.Range("E26").Value = "  -3.86%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  +10.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.16"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.39"
$ws.Range("E30").Value = "  -4.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.93"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.49"
$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.33"
$ws.Range("E33").Value = "  +7.15%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.31"
$ws.Range("E34").Value = "  -0.76%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0785"
$ws.Range("E35").Value = "  -2.60%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").Value = "  -4.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.61"
$ws.Range("E37").Value = "  -5.04%  "

$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.00"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.83"
$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.01"
$ws.Range("E44").Value = "  -3.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0300"
$ws.Range("E45").Value = "  -2.49%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.031.34"
$ws.Range("E46").Value = "  -2.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.74"
$ws.Range("E47").Value = "  -3.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.03"
$ws.Range("E48").Value = "  -3.23%  "

$ws.Range("D49").Value = "2.754.93"
$ws.Range("E49").Value = "  -2.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.78"
$ws.Range("E50").Value = "  -6.71%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.188"
$ws.Range("E51").Value = "  -1.12%  "
